# Add multi-index and multi-table functions into gen datasets script
# Rebuilds the synthetic dataset on Sheet1 (Date, Channel, Metric, Value)
# as a full 7-date x 2-channel x 2-metric table (28 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(46061, "TV", "GRPs", 1),
    @(46089, "TV", "GRPs", 7),
    @(46061, "TV", "Spend", 119),
    @(46082, "TV", "GRPs", 5),
    @(46047, "TV", "Spend", 158),
    @(46068, "TV", "Spend", 130),
    @(46075, "TV", "GRPs", 4),
    @(46082, "Radio", "Spend", 159),
    @(46061, "Radio", "GRPs", 2),
    @(46068, "TV", "GRPs", 1),
    @(46068, "Radio", "GRPs", 1),
    @(46047, "TV", "GRPs", 10),
    @(46054, "TV", "Spend", 159),
    @(46054, "TV", "GRPs", 4),
    @(46047, "Radio", "Spend", 128),
    @(46075, "TV", "Spend", 116),
    @(46082, "Radio", "GRPs", 4),
    @(46047, "Radio", "GRPs", 5),
    @(46089, "Radio", "Spend", 186),
    @(46089, "TV", "Spend", 115),
    @(46075, "Radio", "Spend", 197),
    @(46089, "Radio", "GRPs", 5),
    @(46082, "TV", "Spend", 110),
    @(46054, "Radio", "GRPs", 6),
    @(46061, "Radio", "Spend", 117),
    @(46068, "Radio", "Spend", 147),
    @(46075, "Radio", "GRPs", 9),
    @(46054, "Radio", "Spend", 193)
)

$row = 2
foreach ($rec in $data) {
    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Value = $rec[0]
    $aCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $row = $row + 1
}
